$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab06")

# Update the 8 header cells (row 2, columns C:J): "2012-21" -> "2013-22"
foreach ($col in @("C","D","E","F","G","H","I","J")) {
    $cell = $ws.Range($col + "2")
    $orig = $cell.Value()
    $cell.Value = $orig -replace "2012-21", "2013-22"
}

# Updated data values (row -> column -> new value)
$updates = @{
    "C63" = 97.124107380952395
    "D63" = 96.993711666666698
    "E63" = 97.284585357142902
    "F63" = 0.99524880952381001
    "G63" = 92.493392
    "H63" = 91.049872470588198
    "I63" = 93.952597294117595
    "J63" = 0.96285776470588003

    "C64" = 97.5552858333334
    "D64" = 97.856585416666704
    "E64" = 97.264937083333393
    "F64" = 1.0060150000000001
    "G64" = 92.2380876
    "H64" = 91.667259200000004
    "I64" = 92.8557536
    "J64" = 0.98584000000000005

    "D65" = 93.293322222222301

    "C66" = 90.209619402985098
    "D66" = 89.180483507462696
    "E66" = 91.359260746268703
    "F66" = 0.96591686567164003
    "G66" = 83.185730592592606
    "H66" = 79.986195851851903
    "I66" = 86.552861185185193
    "J66" = 0.90171096296296005

    "C83" = 96.951102121212202
    "D83" = 96.736173333333397
    "E83" = 97.182109393939498
    "F83" = 0.99314181818182001
    "G83" = 91.962249999999997
    "H83" = 90.476356119402993
    "I83" = 93.4991217910448
    "J83" = 0.96030671641790999

    "C87" = 94.573083999999994
    "D87" = 94.368003999999999
    "E87" = 94.762525600000004
    "F87" = 0.99463760000000001
    "G87" = 85.179221200000001
    "H87" = 82.177672000000001
    "I87" = 88.281700400000005
    "J87" = 0.92352319999999999

    "C89" = 98.864615000000001
    "D89" = 98.959063157894803
    "E89" = 98.776679736842098
    "F89" = 1.0017989473684199
    "G89" = 95.982192368420996
    "H89" = 95.313946315789494
    "I89" = 96.665577368421097
    "J89" = 0.98567657894737004

    "G91" = 57.541360322580701

    "C94" = 96.3071985714286
    "D94" = 96.683132857142894
    "E94" = 95.948515
    "F94" = 1.00761428571429
    "G94" = 91.108957333333393
    "H94" = 90.690535999999994
    "I94" = 91.531800666666697
    "J94" = 0.98721800000000004
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
